# Auto-generated edit script: adds LeetCode problems 71-100 (rows 76-106)
# to the tracking sheet, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Write plain cell values for columns A-H (I handled via hyperlinks)
# ---------------------------------------------------------------
# Row 76
$ws.Range("A76").Value = 44348
$ws.Range("B76").Value = 71
$ws.Range("C76").Value = "Simplify Path"
$ws.Range("E76").Value = "stack-pq"
# Row 77
$ws.Range("B77").Value = 72
$ws.Range("C77").Value = "Edit Distance"
$ws.Range("E77").Value = "dp"
$ws.Range("F77").Value = "not yet"
# Row 78
$ws.Range("B78").Value = 73
$ws.Range("C78").Value = "Set Matrix Zeroes"
$ws.Range("E78").Value = "matrix"
# Row 79
$ws.Range("B79").Value = 74
$ws.Range("C79").Value = "Search a 2D Matrix"
$ws.Range("E79").Value = "matrix"
# Row 80
$ws.Range("B80").Value = 75
$ws.Range("C80").Value = "Sort Colors"
$ws.Range("E80").Value = "array"
$ws.Range("F80").Value = "one pass not yet"
# Row 81
$ws.Range("A81").Value = 44349
$ws.Range("B81").Value = 76
$ws.Range("C81").Value = "Minimum Window Substring"
$ws.Range("E81").Value = "string"
$ws.Range("F81").Value = "not yet"
# Row 82
$ws.Range("B82").Value = 77
$ws.Range("C82").Value = "Combinations"
$ws.Range("E82").Value = "backtrack"
$ws.Range("F82").Value = "great problem, transition between iterative and recursive"
# Row 83
$ws.Range("B83").Value = 78
$ws.Range("C83").Value = "Subsets"
$ws.Range("E83").Value = "backtrack"
# Row 84
$ws.Range("B84").Value = 79
$ws.Range("C84").Value = "Word Search"
$ws.Range("E84").Value = "array"
# Row 85
$ws.Range("B85").Value = 80
$ws.Range("C85").Value = "Remove Duplicates from Sorted Array II"
$ws.Range("E85").Value = "array"
$ws.Range("F85").Value = "not yet"
# Row 86
$ws.Range("A86").Value = 44350
$ws.Range("B86").Value = 81
$ws.Range("C86").Value = "Search in Rotated Sorted Array II"
$ws.Range("E86").Value = "binarySearch"
$ws.Range("F86").Value = "not yet"
# Row 87
$ws.Range("B87").Value = 82
$ws.Range("C87").Value = "Remove Duplicates from Sorted List II"
$ws.Range("E87").Value = "LinkedList"
# Row 88
$ws.Range("B88").Value = 83
$ws.Range("C88").Value = "Remove Duplicates from Sorted List"
$ws.Range("E88").Value = "LinkedList"
# Row 89
$ws.Range("B89").Value = 84
$ws.Range("C89").Value = "Largest Rectangle in Histogram"
$ws.Range("E89").Value = "stack-pq"
$ws.Range("F89").Value = "one edge case emitted"
# Row 90
$ws.Range("B90").Value = 85
$ws.Range("C90").Value = "Maximal Rectangle"
$ws.Range("E90").Value = "dp"
$ws.Range("F90").Value = "cannot understand the dp solution"
# Row 91
$ws.Range("A91").Value = 44351
$ws.Range("B91").Value = 86
$ws.Range("C91").Value = "Partition List"
$ws.Range("E91").Value = "LinkedList"
# Row 92
$ws.Range("B92").Value = 87
$ws.Range("C92").Value = "Scramble String"
$ws.Range("E92").Value = "string"
$ws.Range("F92").Value = "dp not yet"
# Row 93
$ws.Range("B93").Value = 88
$ws.Range("C93").Value = "Merge Sorted Array"
$ws.Range("E93").Value = "array"
$ws.Range("F93").Value = "O(m+n)! 想了好一会儿才想出来"
# Row 94
$ws.Range("B94").Value = 89
$ws.Range("C94").Value = "Gray Code"
$ws.Range("F94").Value = "formula"
# Row 95
$ws.Range("B95").Value = 90
$ws.Range("C95").Value = "Subsets II"
$ws.Range("E95").Value = "backtrack"
$ws.Range("F95").Value = "iterative not yet, in essence, this is the same as problem40"
# Row 96
$ws.Range("A96").Value = 44352
$ws.Range("C96").Value = "break"
# Row 97
$ws.Range("A97").Value = 44353
$ws.Range("B97").Value = 91
$ws.Range("C97").Value = "Decode Ways"
$ws.Range("E97").Value = "dp"
$ws.Range("F97").Value = "别人的思路简单一点点"
# Row 98
$ws.Range("B98").Value = 92
$ws.Range("C98").Value = "Reverse Linked List II"
$ws.Range("E98").Value = "LinkedList"
# Row 99
$ws.Range("B99").Value = 93
$ws.Range("C99").Value = "Restore IP Addresses"
$ws.Range("E99").Value = "backtrack"
# Row 100
$ws.Range("B100").Value = 94
$ws.Range("C100").Value = "Binary Tree Inorder Traversal"
$ws.Range("E100").Value = "tree"
# Row 101
$ws.Range("B101").Value = 95
$ws.Range("C101").Value = "Unique Binary Search Trees II"
$ws.Range("E101").Value = "tree"
$ws.Range("F101").Value = "dp not yet, 可以先想想下一道dp咋做"
# Row 102
$ws.Range("B102").Value = 96
$ws.Range("C102").Value = "Unique Binary Search Trees"
$ws.Range("E102").Value = "tree"
# Row 103
$ws.Range("B103").Value = 97
$ws.Range("C103").Value = "Interleaving String"
$ws.Range("E103").Value = "dp"
$ws.Range("F103").Value = "be careful in dp"
# Row 104
$ws.Range("B104").Value = 98
$ws.Range("C104").Value = "Validate Binary Search Tree"
$ws.Range("E104").Value = "tree"
# Row 105
$ws.Range("A105").Value = 44354
$ws.Range("B105").Value = 99
$ws.Range("C105").Value = "Recover Binary Search Tree"
$ws.Range("E105").Value = "tree"
$ws.Range("F105").Value = "可以再做做，稍微想了一下"
# Row 106
$ws.Range("B106").Value = 100
$ws.Range("C106").Value = "Same Tree"
$ws.Range("E106").Value = "tree"

# ---------------------------------------------------------------
# 2) Add hyperlinks for column I (also sets the cell text to the URL)
# ---------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("I76"), "https://leetcode.com/problems/simplify-path/")
$ws.Hyperlinks.Add($ws.Range("I77"), "https://leetcode.com/problems/edit-distance/")
$ws.Hyperlinks.Add($ws.Range("I78"), "https://leetcode.com/problems/set-matrix-zeroes/")
$ws.Hyperlinks.Add($ws.Range("I79"), "https://leetcode.com/problems/search-a-2d-matrix/")
$ws.Hyperlinks.Add($ws.Range("I80"), "https://leetcode.com/problems/sort-colors/")
$ws.Hyperlinks.Add($ws.Range("I81"), "https://leetcode.com/problems/minimum-window-substring/")
$ws.Hyperlinks.Add($ws.Range("I82"), "https://leetcode.com/problems/combinations/")
$ws.Hyperlinks.Add($ws.Range("I83"), "https://leetcode.com/problems/subsets/")
$ws.Hyperlinks.Add($ws.Range("I84"), "https://leetcode.com/problems/word-search/")
$ws.Hyperlinks.Add($ws.Range("I85"), "https://leetcode.com/problems/remove-duplicates-from-sorted-array-ii/")
$ws.Hyperlinks.Add($ws.Range("I86"), "https://leetcode.com/problems/search-in-rotated-sorted-array-ii/")
$ws.Hyperlinks.Add($ws.Range("I87"), "https://leetcode.com/problems/remove-duplicates-from-sorted-list-ii/")
$ws.Hyperlinks.Add($ws.Range("I88"), "https://leetcode.com/problems/remove-duplicates-from-sorted-list/")
$ws.Hyperlinks.Add($ws.Range("I89"), "https://leetcode.com/problems/largest-rectangle-in-histogram/")
$ws.Hyperlinks.Add($ws.Range("I90"), "https://leetcode.com/problems/maximal-rectangle/")
$ws.Hyperlinks.Add($ws.Range("I91"), "https://leetcode.com/problems/partition-list/submissions/")
$ws.Hyperlinks.Add($ws.Range("I92"), "https://leetcode.com/problems/scramble-string/")
$ws.Hyperlinks.Add($ws.Range("I93"), "https://leetcode.com/problems/merge-sorted-array/")
$ws.Hyperlinks.Add($ws.Range("I94"), "https://leetcode.com/problems/gray-code/")
$ws.Hyperlinks.Add($ws.Range("I95"), "https://leetcode.com/problems/subsets-ii/")
$ws.Hyperlinks.Add($ws.Range("I97"), "https://leetcode.com/problems/decode-ways/submissions/")
$ws.Hyperlinks.Add($ws.Range("I98"), "https://leetcode.com/problems/reverse-linked-list-ii/")
$ws.Hyperlinks.Add($ws.Range("I99"), "https://leetcode.com/problems/restore-ip-addresses/")
$ws.Hyperlinks.Add($ws.Range("I100"), "https://leetcode.com/problems/binary-tree-inorder-traversal/")
$ws.Hyperlinks.Add($ws.Range("I101"), "https://leetcode.com/problems/unique-binary-search-trees-ii/")
$ws.Hyperlinks.Add($ws.Range("I102"), "https://leetcode.com/problems/unique-binary-search-trees/")
$ws.Hyperlinks.Add($ws.Range("I103"), "https://leetcode.com/problems/interleaving-string/")
$ws.Hyperlinks.Add($ws.Range("I104"), "https://leetcode.com/problems/validate-binary-search-tree/")
$ws.Hyperlinks.Add($ws.Range("I105"), "https://leetcode.com/problems/recover-binary-search-tree/")
$ws.Hyperlinks.Add($ws.Range("I106"), "https://leetcode.com/problems/same-tree/")

# ---------------------------------------------------------------
# 3) Fix up hyperlink cell formatting: Hyperlinks.Add creates its own
#    "Hyperlink" style variant; copy the format from an existing,
#    already-correct hyperlink cell (I75) so the new cells reuse the
#    same style index instead of bloating the style table.
# ---------------------------------------------------------------
$ws.Range("I75").Copy()
$ws.Range("I76:I106").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 4) Sheet view / selection / dimension bookkeeping to match target
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 88
$ws.Range("B107").Select()
